# Raul's Log.xlsx - add THURSDAY (Aug 4, 2016) and FRIDAY (Aug 5, 2016)
# sections to the bottom of the "Logs" sheet, following the same layout
# used by the existing day-of-week blocks (e.g. the "WEDNESDAY" block
# immediately above them).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# ---------------------------------------------------------------------
# 1) Clone cell formatting from existing template rows so the new rows
#    look identical to the rest of the log (day headers + data rows).
# ---------------------------------------------------------------------

# Day header template (row 190 = "WEDNESDAY" banner spanning A:F)
$ws.Range("A190:F190").Copy() | Out-Null
$ws.Range("A196:F196").PasteSpecial(-4122) | Out-Null
$ws.Range("A203:F203").PasteSpecial(-4122) | Out-Null

# Data-row template with only columns A:E populated (row 191)
$ws.Range("A191:E191").Copy() | Out-Null
$ws.Range("A197:E197").PasteSpecial(-4122) | Out-Null
$ws.Range("A198:E198").PasteSpecial(-4122) | Out-Null
$ws.Range("A199:E199").PasteSpecial(-4122) | Out-Null
$ws.Range("A200:E200").PasteSpecial(-4122) | Out-Null
$ws.Range("A205:E205").PasteSpecial(-4122) | Out-Null
$ws.Range("A206:E206").PasteSpecial(-4122) | Out-Null

# Data-row template with columns A:F populated (row 184)
$ws.Range("A184:F184").Copy() | Out-Null
$ws.Range("A204:F204").PasteSpecial(-4122) | Out-Null
$ws.Range("A207:F207").PasteSpecial(-4122) | Out-Null
$ws.Range("A208:F208").PasteSpecial(-4122) | Out-Null
$ws.Range("A209:F209").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) THURSDAY block (Aug 4, 2016)
# ---------------------------------------------------------------------

$ws.Range("B196").Value2 = "THURSDAY"

$ws.Range("A197").Value2 = "Crestron Logout"
$ws.Range("B197").Value2 = 42586
$ws.Range("C197").Value2 = "1630"
$ws.Range("D197").Value2 = "OSG"
$ws.Range("E197").Value2 = "1001"

$ws.Range("A198").Value2 = "Crestron Logout"
$ws.Range("B198").Value2 = 42586
$ws.Range("C198").Value2 = "1630"
$ws.Range("D198").Value2 = "OSG"
$ws.Range("E198").Value2 = "1002"

$ws.Range("A199").Value2 = "Crestron Logout"
$ws.Range("B199").Value2 = 42586
$ws.Range("C199").Value2 = "1630"
$ws.Range("D199").Value2 = "OSG"
$ws.Range("E199").Value2 = "2001"

$ws.Range("A200").Value2 = "Crestron Logout"
$ws.Range("B200").Value2 = 42586
$ws.Range("C200").Value2 = "1630"
$ws.Range("D200").Value2 = "OSG"
$ws.Range("E200").Value2 = "2003"

# ---------------------------------------------------------------------
# 3) FRIDAY block (Aug 5, 2016)
# ---------------------------------------------------------------------

$ws.Range("B203").Value2 = "FRIDAY"

$ws.Range("A204").Value2 = "Pickup Skype Kit"
$ws.Range("B204").Value2 = 42587
$ws.Range("C204").Value2 = "1730"
$ws.Range("D204").Value2 = "ACE"
$ws.Range("E204").Value2 = "003"
$ws.Range("F204").Value2 = "Return camera and tripod to ACE 015 storeroom."

$ws.Range("A205").Value2 = "Crestron Logout"
$ws.Range("B205").Value2 = 42587
$ws.Range("C205").Value2 = "1630"
$ws.Range("D205").Value2 = "OSG"
$ws.Range("E205").Value2 = "1003"

$ws.Range("A206").Value2 = "Crestron Logout"
$ws.Range("B206").Value2 = 42587
$ws.Range("C206").Value2 = "1630"
$ws.Range("D206").Value2 = "OSG"
$ws.Range("E206").Value2 = "2002"

$ws.Range("A207").Value2 = "Pickup Mic"
$ws.Range("B207").Value2 = 42587
$ws.Range("C207").Value2 = "2100"
$ws.Range("D207").Value2 = "OSG"
$ws.Range("E207").Value2 = "1005"
$ws.Range("F207").Value2 = "Pick up Lecturn mic and stand and cable from Osgoode Moot Court Room. Put in closet just to right of PC cabinet behind podium."
$ws.Rows.Item(207).RowHeight = 45

$ws.Range("A208").Value2 = "Other"
$ws.Range("B208").Value2 = 42587
$ws.Range("C208").Value2 = "2100"
$ws.Range("D208").Value2 = "OSG"
$ws.Range("E208").Value2 = "1005"
$ws.Range("F208").Value2 = "Return neck mic to cabinet drawer and shut drawer."

$ws.Range("A209").Value2 = "Crestron Logout"
$ws.Range("B209").Value2 = 42587
$ws.Range("C209").Value2 = "2100"
$ws.Range("D209").Value2 = "OSG"
$ws.Range("E209").Value2 = "1005"
$ws.Range("F209").Value2 = "Log off crestron in Osgoode Moot Court."

# ---------------------------------------------------------------------
# 4) Restore the view state (scroll position / active cell) the author
#    left the sheet in after typing the new rows.
# ---------------------------------------------------------------------

$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 181
$win.ScrollColumn = 1
$ws.Range("F212").Select() | Out-Null
